# "User Screen Code Updated"
#
# Semantic changes in this commit:
#  - AddUser!A5: "User_TC004" -> "User_TC004A"
#  - UsersGrid!J12: "Recod not found" (typo) -> "Record not found"
#  - UsersGrid!J14: "Recod not found" (typo) -> "Record not found"
#  - Active sheet switches from "UsersGrid" to "AddUser" (first tab becomes
#    the active one again), with the on-sheet selections updated too:
#      AddUser   selection: C4  -> A5
#      UsersGrid selection: A14 -> J6

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AddUser")
$ws2 = $wb.Worksheets.Item("UsersGrid")

# --- cell value corrections ---
$ws1.Range("A5").Value  = "User_TC004A"
$ws2.Range("J12").Value = "Record not found"
$ws2.Range("J14").Value = "Record not found"

# --- selection / active-sheet updates ---
# Set UsersGrid's remembered selection first (it is being deactivated).
$ws2.Activate()
$ws2.Range("J6").Select()

# Finish on AddUser so it becomes the workbook's active tab, with its own
# selection restored to A5.
$ws1.Activate()
$ws1.Range("A5").Select()
